# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp
# - A handful of countries changed rank (sorted desc by "Casos totales") now that
#   their case counts were refreshed, so the country label in some rows needs to
#   swap with its neighbour(s)
# - Apply the refreshed Covid figures (Casos totales/Nuevos casos/Casos
#   activos/Recuperados/Casos criticos/Muertes hoy/Muertes) for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 12:47"

# Country name re-rankings (rows whose position swapped due to updated case counts)
$ws.Range("A20").Value = "Banglades"
$ws.Range("A21").Value = "Canada"
$ws.Range("A33").Value = "Indonesia"
$ws.Range("A34").Value = "Singapur"
$ws.Range("A85").Value = "El Salvador"
$ws.Range("A86").Value = "Luxemburgo"
$ws.Range("A87").Value = "Hungria"
$ws.Range("A145").Value = "Estado de Palestina"
$ws.Range("A146").Value = "Malaui"
$ws.Range("A147").Value = "Benin"
$ws.Range("A148").Value = "Suazilandia"
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"

# Updated numeric figures (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# Row 4: B4=2234854, C4=383, E4=1196115, G4=2, H4=119943
$ws.Range("B4").Value = 2234854
$ws.Range("C4").Value = 383
$ws.Range("E4").Value = 1196115
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 119943

# Row 12: B12=197647, C12=2596, D12=156991, E12=31384, G12=87, H12=9272
$ws.Range("B12").Value = 197647
$ws.Range("C12").Value = 2596
$ws.Range("D12").Value = 156991
$ws.Range("E12").Value = 31384
$ws.Range("G12").Value = 87
$ws.Range("H12").Value = 9272

# Row 20: B20=102292, C20=3803, D20=40164, E20=60785, G20=38, H20=1343
$ws.Range("B20").Value = 102292
$ws.Range("C20").Value = 3803
$ws.Range("D20").Value = 40164
$ws.Range("E20").Value = 60785
$ws.Range("G20").Value = 38
$ws.Range("H20").Value = 1343

# Row 21: B21=99853, D21=62017, E21=29582, H21=8254
$ws.Range("B21").Value = 99853
$ws.Range("D21").Value = 62017
$ws.Range("E21").Value = 29582
$ws.Range("H21").Value = 8254

# Row 25: B25=60348, C25=104, D25=16724, E25=33941, G25=8, H25=9683
$ws.Range("B25").Value = 60348
$ws.Range("C25").Value = 104
$ws.Range("D25").Value = 16724
$ws.Range("E25").Value = 33941
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 9683

# Row 33: B33=42762, C33=1331, D33=16798, E33=23625, G33=63, H33=2339
$ws.Range("B33").Value = 42762
$ws.Range("C33").Value = 1331
$ws.Range("D33").Value = 16798
$ws.Range("E33").Value = 23625
$ws.Range("G33").Value = 63
$ws.Range("H33").Value = 2339

# Row 34: B34=41473, C34=257, D34=31938, E34=9509, H34=26
$ws.Range("B34").Value = 41473
$ws.Range("C34").Value = 257
$ws.Range("D34").Value = 31938
$ws.Range("E34").Value = 9509
$ws.Range("H34").Value = 26

# Row 39: B39=31200, C39=13, E39=344
$ws.Range("B39").Value = 31200
$ws.Range("C39").Value = 13
$ws.Range("E39").Value = 344

# Row 47: B47=23080, C47=320, D47=16308, E47=5299, G47=22, H47=1473
$ws.Range("B47").Value = 23080
$ws.Range("C47").Value = 320
$ws.Range("D47").Value = 16308
$ws.Range("E47").Value = 5299
$ws.Range("G47").Value = 22
$ws.Range("H47").Value = 1473

# Row 50: E50=5724, G50=3, H50=52
$ws.Range("E50").Value = 5724
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 52

# Row 68: B68=9042, C68=45, D68=7999, E68=830
$ws.Range("B68").Value = 9042
$ws.Range("C68").Value = 45
$ws.Range("D68").Value = 7999
$ws.Range("E68").Value = 830

# Row 85: B85=4200, C85=134, D85=2235, E85=1883, G85=3, H85=82
$ws.Range("B85").Value = 4200
$ws.Range("C85").Value = 134
$ws.Range("D85").Value = 2235
$ws.Range("E85").Value = 1883
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 82

# Row 86: B86=4085, C86=0, D86=3935, E86=40, G86=0, H86=110
$ws.Range("B86").Value = 4085
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 3935
$ws.Range("E86").Value = 40
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 110

# Row 87: B87=4079, C87=1, D87=2564, E87=947, G87=1, H87=568
$ws.Range("B87").Value = 4079
$ws.Range("C87").Value = 1
$ws.Range("D87").Value = 2564
$ws.Range("E87").Value = 947
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 568

# Row 145: B145=579, C145=24, D145=415, E145=161, H145=3
$ws.Range("B145").Value = 579
$ws.Range("C145").Value = 24
$ws.Range("D145").Value = 415
$ws.Range("E145").Value = 161
$ws.Range("H145").Value = 3

# Row 146: D146=73, E146=493, H146=6
$ws.Range("D146").Value = 73
$ws.Range("E146").Value = 493
$ws.Range("H146").Value = 6

# Row 147: B147=572, D147=237, E147=326, H147=9
$ws.Range("B147").Value = 572
$ws.Range("D147").Value = 237
$ws.Range("E147").Value = 326
$ws.Range("H147").Value = 9

# Row 148: B148=563, C148=0, D148=262, E148=297, H148=4
$ws.Range("B148").Value = 563
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 262
$ws.Range("E148").Value = 297
$ws.Range("H148").Value = 4

# Row 190: B190=39, C190=3, E190=20
$ws.Range("B190").Value = 39
$ws.Range("C190").Value = 3
$ws.Range("E190").Value = 20

# Row 208: D208=12, H208=0
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209: D209=11, H209=1
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
